# Scheduled-runner refresh of recalculated Leve-profit figures (currentAveragePrice*,
# LevePrice*/LeveProfit* columns H..N) across the per-job sheets. Values below are the
# updated (refreshed) market-board-derived numbers; no structural/formula changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 12995.546
$ws.Range("J88").Value = 16619.625
$ws.Range("L88").Value = 16619.625
$ws.Range("N88").Value = -17431.625
$ws.Range("H91").Value = 12995.546
$ws.Range("J91").Value = 16619.625
$ws.Range("L91").Value = 16619.625
$ws.Range("N91").Value = -19427.625
$ws.Range("H113").Value = 5749.75
$ws.Range("I113").Value = 5749.75
$ws.Range("K113").Value = 5749.75
$ws.Range("M113").Value = -2495.75
$ws.Range("H132").Value = 2161.7856
$ws.Range("I132").Value = 2161.7856
$ws.Range("K132").Value = 6485.3568
$ws.Range("M132").Value = -3955.3568
$ws.Range("H135").Value = 659.1818
$ws.Range("I135").Value = 620.3
$ws.Range("J135").Value = 1048
$ws.Range("K135").Value = 5582.7
$ws.Range("L135").Value = 9432
$ws.Range("M135").Value = -3047.7
$ws.Range("N135").Value = -14502
$ws.Range("H138").Value = 2856.625
$ws.Range("J138").Value = 3198.9565
$ws.Range("L138").Value = 9596.869499999999
$ws.Range("N138").Value = -19876.8695
$ws.Range("H141").Value = 5382.4165
$ws.Range("I141").Value = 917.8
$ws.Range("J141").Value = 8571.429
$ws.Range("K141").Value = 2753.4
$ws.Range("L141").Value = 25714.287
$ws.Range("M141").Value = 2426.6
$ws.Range("N141").Value = -36074.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4522.927
$ws.Range("I32").Value = 3850.2285
$ws.Range("J32").Value = 8447
$ws.Range("K32").Value = 3850.2285
$ws.Range("L32").Value = 8447
$ws.Range("M32").Value = -3563.2285
$ws.Range("N32").Value = -9021
$ws.Range("H132").Value = 5753.875
$ws.Range("I132").Value = 3859.5715
$ws.Range("K132").Value = 11578.7145
$ws.Range("M132").Value = -9048.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 44060.4
$ws.Range("J103").Value = 44060.4
$ws.Range("L103").Value = 44060.4
$ws.Range("N103").Value = -46404.4
$ws.Range("H134").Value = 2849.6667
$ws.Range("I134").Value = 1420
$ws.Range("J134").Value = 9998
$ws.Range("K134").Value = 4260
$ws.Range("L134").Value = 29994
$ws.Range("M134").Value = -1725
$ws.Range("N134").Value = -35064

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24542.06
$ws.Range("I31").Value = 2673.7297
$ws.Range("K31").Value = 2673.7297
$ws.Range("M31").Value = -2378.7297
$ws.Range("H34").Value = 24542.06
$ws.Range("I34").Value = 2673.7297
$ws.Range("K34").Value = 2673.7297
$ws.Range("M34").Value = -2471.7297
$ws.Range("H62").Value = 10747.75
$ws.Range("I62").Value = 4985
$ws.Range("J62").Value = 12668.667
$ws.Range("K62").Value = 4985
$ws.Range("L62").Value = 12668.667
$ws.Range("M62").Value = -4361
$ws.Range("N62").Value = -13916.667
$ws.Range("H65").Value = 10747.75
$ws.Range("I65").Value = 4985
$ws.Range("J65").Value = 12668.667
$ws.Range("K65").Value = 24925
$ws.Range("L65").Value = 63343.335
$ws.Range("M65").Value = -21805
$ws.Range("N65").Value = -69583.33499999999
$ws.Range("H74").Value = 95517
$ws.Range("J74").Value = 184842
$ws.Range("L74").Value = 184842
$ws.Range("N74").Value = -186590
$ws.Range("H77").Value = 95517
$ws.Range("J77").Value = 184842
$ws.Range("L77").Value = 554526
$ws.Range("N77").Value = -563262
$ws.Range("H86").Value = 6374.7144
$ws.Range("I86").Value = 4664.5454
$ws.Range("J86").Value = 8255.9
$ws.Range("K86").Value = 4664.5454
$ws.Range("L86").Value = 8255.9
$ws.Range("M86").Value = -3541.5454
$ws.Range("N86").Value = -10501.9
$ws.Range("H89").Value = 6374.7144
$ws.Range("I89").Value = 4664.5454
$ws.Range("J89").Value = 8255.9
$ws.Range("K89").Value = 23322.727
$ws.Range("L89").Value = 41279.5
$ws.Range("M89").Value = -17706.727
$ws.Range("N89").Value = -52511.5
$ws.Range("H134").Value = 8014.75
$ws.Range("I134").Value = 5225
$ws.Range("K134").Value = 15675
$ws.Range("M134").Value = -13140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2744.7778
$ws.Range("H129").Value = 7580070.5
$ws.Range("I129").Value = 606
$ws.Range("J129").Value = 11911193
$ws.Range("K129").Value = 1818
$ws.Range("L129").Value = 35733579
$ws.Range("M129").Value = 3182
$ws.Range("N129").Value = -35743579

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9475.375
$ws.Range("I80").Value = 5949.25
$ws.Range("J80").Value = 13001.5
$ws.Range("K80").Value = 5949.25
$ws.Range("L80").Value = 13001.5
$ws.Range("M80").Value = -4951.25
$ws.Range("N80").Value = -14997.5
$ws.Range("H83").Value = 9475.375
$ws.Range("I83").Value = 5949.25
$ws.Range("J83").Value = 13001.5
$ws.Range("K83").Value = 29746.25
$ws.Range("L83").Value = 65007.5
$ws.Range("M83").Value = -24754.25
$ws.Range("N83").Value = -74991.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1668650
$ws.Range("J55").Value = 2564.2307
$ws.Range("L55").Value = 2564.2307
$ws.Range("N55").Value = -2910.2307
$ws.Range("H68").Value = 5221.5
$ws.Range("J68").Value = 16334
$ws.Range("L68").Value = 16334
$ws.Range("N68").Value = -17832
$ws.Range("H71").Value = 5221.5
$ws.Range("J71").Value = 16334
$ws.Range("L71").Value = 81670
$ws.Range("N71").Value = -89158
$ws.Range("H93").Value = 2476.1667
$ws.Range("I93").Value = 1870.2
$ws.Range("J93").Value = 3486.111
$ws.Range("K93").Value = 1870.2
$ws.Range("L93").Value = 3486.111
$ws.Range("M93").Value = -622.2
$ws.Range("N93").Value = -5982.111
$ws.Range("H103").Value = 20097.5
$ws.Range("I103").Value = 50000
$ws.Range("K103").Value = 50000
$ws.Range("M103").Value = -48828
$ws.Range("H108").Value = 84967
$ws.Range("J108").Value = 84967
$ws.Range("L108").Value = 84967
$ws.Range("N108").Value = -92647
$ws.Range("H136").Value = 6550.136
$ws.Range("I136").Value = 2922.0908
$ws.Range("J136").Value = 10178.182
$ws.Range("K136").Value = 8766.2724
$ws.Range("L136").Value = 30534.546
$ws.Range("M136").Value = -6216.2724
$ws.Range("N136").Value = -35634.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
